$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts Trading/Salary/Freelancing down by one)
$ws.Rows.Item(2).Insert()

# Populate the new "Teacher" income row
$ws.Range("A2").Value = "Teacher"
$ws.Range("B2").Value = 50000
$ws.Range("C2").Value = 45755.22928240741

# The insert leaves the new row's date cell with General formatting; copy the
# date number format down from the row below (now "Trading", originally row 2)
# so C2 matches the other date cells in column C.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
